$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unmerge the old title-bar merge and clear its centered-alignment style ---
$ws.Range("A1:E1").UnMerge()
$ws.Range("A1:E1").ClearFormats()
$ws.Range("E1").ClearContents()

# --- Seed unique text values in the same order the original author typed them,
#     so shared-string indices line up with the authored workbook. ---
$ws.Range("A1").Value = "Selected clothing"
$ws.Range("A2").Value = "Jacket"
$ws.Range("A5").Value = "Pants"
$ws.Range("A8").Value = "Scarf"
$ws.Range("B1").Value = "Variable 1 (buttons, pockets, threads) "
$ws.Range("D1").Value = "Win count"
$ws.Range("C1").Value = "Variable 2 (thickness)"
$ws.Range("D2").Value = "full win"
$ws.Range("D3").Value = "half win"
$ws.Range("D4").Value = "loss"

# --- Fill remaining "Selected clothing" labels ---
$ws.Range("A3").Value = "Jacket"
$ws.Range("A4").Value = "Jacket"
$ws.Range("A6").Value = "Pants"
$ws.Range("A7").Value = "Pants"
$ws.Range("A9").Value = "Scarf"
$ws.Range("A10").Value = "Scarf"

# --- Fill remaining win-count labels ---
$ws.Range("D5").Value = "full win"
$ws.Range("D6").Value = "half win"
$ws.Range("D7").Value = "loss"
$ws.Range("D8").Value = "full win"
$ws.Range("D9").Value = "half win"
$ws.Range("D10").Value = "loss"

# --- Variable 1 (buttons, pockets, threads) numeric column ---
$ws.Range("B2").Value = 4
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 7
$ws.Range("B5").Value = 2
$ws.Range("B6").Value = 4
$ws.Range("B7").Value = 10
$ws.Range("B8").Value = 15
$ws.Range("B9").Value = 8
$ws.Range("B10").Value = 9

# --- Variable 2 (thickness) numeric column ---
$ws.Range("C2").Value = 1.5
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 2.5
$ws.Range("C5").Value = 0.1
$ws.Range("C6").Value = 1
$ws.Range("C7").Value = 10
$ws.Range("C8").Value = 1
$ws.Range("C9").Value = 2
$ws.Range("C10").Value = 0.9999999

# --- Column widths (chars) -- values chosen so the exported OOXML <col width>
# (which this runtime quantizes to 1/6-character pixel steps) lands as close as
# possible to the authored widths of 20.83203125 / 32.6640625 / 20.83203125 ---
$ws.Columns.Item(1).ColumnWidth = 20
$ws.Columns.Item(2).ColumnWidth = 31.8333333333
$ws.Range($ws.Cells.Item(1,3), $ws.Cells.Item(1,15)).EntireColumn.ColumnWidth = 20

# --- Selection ---
$ws.Range("E10").Select()
